$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("M8").Value2 = '{"success":true, "data":[ {"id":"SPECIES","label":"Informations sur les esp\u00e8ces","definition":"Informations sur les donn\u00e9es des esp\u00e8ces","is_default":"1"},{"id":"TREES","label":"D\u00e9tails sur les arbres","definition":"D\u00e9tails sur les donn\u00e9es de chaque arbre","is_default":"0"}]}'

$ws.Range("M10").Value2 = '{"success":true,"total":2, "root":[{"inputType":"TEXT","isCriteria":"1","isResult":"1","isDefaultCriteria":"0","isDefaultResult":"1","defaultValue":null,"decimals":null,"mask":null,"editable":null,"insertable":null,"isPK":null,"required":null,"position":2,"data":"PLOT_CODE","format":"PLOT_FORM","label":"Plot code","unit":"PLOT_CODE","type":"STRING","subtype":null,"definition":"The identifier of the plot for the country","value":null,"valueLabel":null},{"inputType":"CHECKBOX","isCriteria":"1","isResult":"1","isDefaultCriteria":"1","isDefaultResult":"1","defaultValue":"1","decimals":null,"mask":null,"editable":null,"insertable":null,"isPK":null,"required":null,"position":5,"data":"IS_FOREST_PLOT","format":"PLOT_FORM","label":"Is a forest plot","unit":"IS_FOREST_PLOT","type":"STRING","subtype":null,"definition":"True if the plot is a forest plot","value":null,"valueLabel":null}]}'

$ws.Range("M17").Value2 = '{"success":true, "data":[{"text":"Habitats littoraux et halophiles","id":"1"},{"text":"Milieux aquatiques non marins","id":"2"},{"text":"Landes, frutic\u00e9es, pelouses et prairies","id":"3"},{"text":"For\u00eats","id":"4"},{"text":"Tourbi\u00e8res et marais","id":"5"},{"text":"Rochers continentaux, \u00e9boulis et sables","id":"6"},{"text":"Terres agricoles et paysages artificiels","id":"8"}]}'

$ws.Range("M19").Value2 = '{"success":true, "data":[{"text":"Vivant","id":"349525","vernacularName":null,"isReference":"1"}]}'

$ws.Range("M21").Value2 = '{"success":true,"formats":[{"title":"Localisation","children_count":1,"id":"SCHEMA\/RAW_DATA\/FORMAT\/LOCATION_DATA\/PROVIDER_ID\/1\/PLOT_CODE\/01575-14060-4-0T","fields":[{"label":"Commune","value":["CORMORANCHE-SUR-SAONE (01123)","GARNERANS (01167)"],"inputType":"SELECT","type":"ARRAY"},{"label":"D\u00e9partement","value":"AIN","inputType":"SELECT","type":"CODE"}],"editURL":"SCHEMA\/RAW_DATA\/FORMAT\/LOCATION_DATA\/PROVIDER_ID\/1\/PLOT_CODE\/01575-14060-4-0T"},{"title":"Donn\u00e9es du point","children_count":0,"id":"SCHEMA\/RAW_DATA\/FORMAT\/PLOT_DATA\/PROVIDER_ID\/1\/PLOT_CODE\/01575-14060-4-0T\/CYCLE\/5","fields":[{"label":"Identifiant du fournisseur de donn\u00e9es","value":"France","inputType":"SELECT","type":"CODE"},{"label":"Code du point","value":"01575-14060-4-0T","inputType":"TEXT","type":"STRING"},{"label":"Cycle","value":"5","inputType":"TEXT","type":"STRING"},{"label":"Date d''inventaire","value":"2007\/11\/27","inputType":"DATE","type":"DATE"},{"label":"Est un point for\u00eat","value":"1","inputType":"CHECKBOX","type":"STRING"},{"label":"Corine biotope","value":["Habitats littoraux et halophiles"],"inputType":"TREE","type":"ARRAY"},{"label":"Commentaire","value":"test data plot","inputType":"TEXT","type":"STRING"}],"editURL":"SCHEMA\/RAW_DATA\/FORMAT\/PLOT_DATA\/PROVIDER_ID\/1\/PLOT_CODE\/01575-14060-4-0T\/CYCLE\/5"}],"title":"Donn\u00e9es du point (1_01575-14060-4-0T_5)","maps1":{"title":"image","urls":[{"url":"http:\/\/wxs-i.ign.fr\/7gr31kqe5xttprd2g7zbkqgo\/geoportail\/r\/wms?LAYERS=ORTHOIMAGERY.ORTHOPHOTOS&amp;TRANSPARENT=true&amp;FORMAT=image%2Fpng&amp;SERVICE=WMS&amp;VERSION=1.3.0&amp;REQUEST=GetMap&amp;STYLES=&amp;BBOX=531065.24180649,5812149.7543612,541065.24180649,5822149.7543612&amp;WIDTH=300&amp;HEIGHT=300&amp;map.scalebar=STATUS+embed&amp;SESSION_ID=mlk6a2jgaij2vmqtv0n4njcruvaoj7qj&amp;PROVIDER_ID=1&amp;PLOT_CODE=01575-14060-4-0T&amp;CRS=EPSG%3A3857"},{"url":"http:\/\/localhost\/mapserv-ogam?LAYERS=result_location&amp;TRANSPARENT=true&amp;FORMAT=image%2Fpng&amp;SERVICE=WMS&amp;VERSION=1.1.1&amp;REQUEST=GetMap&amp;STYLES=&amp;BBOX=531065.24180649,5812149.7543612,541065.24180649,5822149.7543612&amp;WIDTH=300&amp;HEIGHT=300&amp;map.scalebar=STATUS+embed&amp;SESSION_ID=mlk6a2jgaij2vmqtv0n4njcruvaoj7qj&amp;PROVIDER_ID=1&amp;PLOT_CODE=01575-14060-4-0T&amp;SRS=EPSG%3A3857"}]},"maps2":{"title":"overview","urls":[{"url":"http:\/\/wxs-i.ign.fr\/7gr31kqe5xttprd2g7zbkqgo\/geoportail\/r\/wms?LAYERS=ORTHOIMAGERY.ORTHOPHOTOS&amp;TRANSPARENT=true&amp;FORMAT=image%2Fpng&amp;SERVICE=WMS&amp;VERSION=1.3.0&amp;REQUEST=GetMap&amp;STYLES=&amp;BBOX=436065.24180649,5717149.7543612,636065.24180649,5917149.7543612&amp;WIDTH=300&amp;HEIGHT=300&amp;map.scalebar=STATUS+embed&amp;SESSION_ID=mlk6a2jgaij2vmqtv0n4njcruvaoj7qj&amp;PROVIDER_ID=1&amp;PLOT_CODE=01575-14060-4-0T&amp;CRS=EPSG%3A3857"},{"url":"http:\/\/www.ogam.fr\/mapProxy.php?LAYERS=nuts_0&amp;TRANSPARENT=true&amp;FORMAT=image%2Fpng&amp;SERVICE=WMS&amp;VERSION=1.1.1&amp;REQUEST=GetMap&amp;STYLES=&amp;BBOX=436065.24180649,5717149.7543612,636065.24180649,5917149.7543612&amp;WIDTH=300&amp;HEIGHT=300&amp;map.scalebar=STATUS+embed&amp;SESSION_ID=mlk6a2jgaij2vmqtv0n4njcruvaoj7qj&amp;PROVIDER_ID=1&amp;PLOT_CODE=01575-14060-4-0T&amp;SRS=EPSG%3A3857"},{"url":"http:\/\/localhost\/mapserv-ogam?LAYERS=plot_location_detail&amp;TRANSPARENT=true&amp;FORMAT=image%2Fpng&amp;SERVICE=WMS&amp;VERSION=1.1.1&amp;REQUEST=GetMap&amp;STYLES=&amp;BBOX=436065.24180649,5717149.7543612,636065.24180649,5917149.7543612&amp;WIDTH=300&amp;HEIGHT=300&amp;map.scalebar=STATUS+embed&amp;SESSION_ID=mlk6a2jgaij2vmqtv0n4njcruvaoj7qj&amp;PROVIDER_ID=1&amp;PLOT_CODE=01575-14060-4-0T&amp;SRS=EPSG%3A3857"}]}}'

$ws.Range("M28").Value2 = '{"success":true,"services":[{"name":"geoportal_wmts", "config":{"urls": ["http://wxs-i.ign.fr/7gr31kqe5xttprd2g7zbkqgo/geoportail/wmts?"],"params":{"SERVICE":"WMTS","VERSION":"1.0.0","REQUEST":"getTile","style":"normal","matrixSet":"PM","requestEncoding":"KVP","maxExtent":[-20037508, -20037508, 20037508, 20037508],"serverResolutions":[156543.033928,78271.516964,39135.758482,19567.879241,9783.939621,4891.969810,2445.984905,1222.992453,611.496226,305.748113,152.874057,76.437028,38.218514,19.109257,9.554629,4.777302,2.388657,1.194329,0.597164,0.298582,0.149291,0.074646],"tileOrigin":[-20037508,20037508]}}},{"name":"mapProxy_wfs", "config":{"urls":["http://www.ogam.fr/mapProxy.php?"],"params":{"SERVICE":"WFS","VERSION":"1.1.0","REQUEST":"GetFeature"}}},{"name":"geoportal_wms", "config":{"urls":["http://wxs-i.ign.fr/7gr31kqe5xttprd2g7zbkqgo/geoportail/r/wms?"],"params":{"SERVICE":"WMS","VERSION":"1.3.0","REQUEST":"GetMap"}}},{"name":"legend_mapProxy", "config":{"urls":["http://www.ogam.fr/mapProxy.php?"],"params":{"SERVICE":"WMS","VERSION":"1.1.1","REQUEST":"GetLegendGraphic"}}},{"name":"proxy_wfs", "config":{"urls":["http://www.ogam.fr/proxy/getwfs?"],"params":{"SERVICE":"WFS","VERSION":"1.0.0","REQUEST":"GetFeature"}}},{"name":"local_mapserver", "config":{"urls":["http://localhost/mapserv-ogam?"],"params":{"SERVICE":"WMS","VERSION":"1.1.1","REQUEST":"GetMap"}}},{"name":"mapProxy", "config":{"urls":["http://www.ogam.fr/mapProxy.php?"],"params":{"SERVICE":"WMS","VERSION":"1.1.1","REQUEST":"GetMap"}}}],"layers":[{"singleTile":false, "name":"nuts_0", "viewServiceName":"mapProxy", "featureServiceName":"mapProxy_wfs", "legendServiceName":"legend_mapProxy", "params":{"layers" : ["nuts_0"], "transparent": true, "format": "image/PNG", "isHidden": false, "isDisabled": false, "isChecked": true, "activateType": "NONE", "hasSLD": false, "session_id": "mlk6a2jgaij2vmqtv0n4njcruvaoj7qj", "provider_id": "1"}, "options":{"buffer": 0, "nodeGroup": "6", "transitionEffect": "resize", "visibility": true, "isBaseLayer": false, "opacity": 1,"label":"Frontières", "resolutions": [8819.439681947, 3527.7758727788, 1763.8879363894, 881.9439681947, 352.77758727788, 176.38879363894]}},{"singleTile":false, "name":"plot_location_detail", "viewServiceName":"mapProxy", "featureServiceName":"", "legendServiceName":"", "params":{"layers" : ["plot_location_detail"], "transparent": true, "format": "image/PNG", "isHidden": true, "isDisabled": false, "isChecked": false, "activateType": "NONE", "hasSLD": false, "session_id": "mlk6a2jgaij2vmqtv0n4njcruvaoj7qj", "provider_id": "1"}, "options":{"buffer": 0, "nodeGroup": "-1", "visibility": false, "isBaseLayer": false, "opacity": 1,"label":"Détails sur les localisations"}},{"singleTile":true, "name":"all_locations", "viewServiceName":"mapProxy", "featureServiceName":"mapProxy_wfs", "legendServiceName":"legend_mapProxy", "params":{"layers" : ["all_locations"], "transparent": true, "format": "image/PNG", "isHidden": false, "isDisabled": false, "isChecked": true, "activateType": "NONE", "hasSLD": false, "session_id": "mlk6a2jgaij2vmqtv0n4njcruvaoj7qj", "provider_id": "1"}, "options":{"buffer": 0, "nodeGroup": "-1", "visibility": true, "isBaseLayer": false, "opacity": 1,"label":"Localisations"}},{"singleTile":false, "name":"ORTHOIMAGERY.ORTHOPHOTOS", "viewServiceName":"geoportal_wmts", "featureServiceName":"", "legendServiceName":"", "params":{"layers" : ["ORTHOIMAGERY.ORTHOPHOTOS"], "transparent": true, "format": "image/jpeg", "isHidden": false, "isDisabled": false, "isChecked": false, "activateType": "NONE", "hasSLD": false, "session_id": "mlk6a2jgaij2vmqtv0n4njcruvaoj7qj", "provider_id": "1"}, "options":{"buffer": 0, "nodeGroup": "-1", "transitionEffect": "resize", "visibility": false, "isBaseLayer": false, "opacity": "0.5","label":"BD Ortho"}},{"singleTile":true, "name":"result_location", "viewServiceName":"mapProxy", "featureServiceName":"", "legendServiceName":"legend_mapProxy", "params":{"layers" : ["result_location"], "transparent": true, "format": "image/PNG", "isHidden": false, "isDisabled": true, "isChecked": false, "activateType": "REQUEST", "hasSLD": false, "session_id": "mlk6a2jgaij2vmqtv0n4njcruvaoj7qj", "provider_id": "1"}, "options":{"buffer": 0, "nodeGroup": "-1", "visibility": false, "isBaseLayer": false, "opacity": 1,"label":"Results"}}]}'

$ws.Range("M29").Value2 = '{"success":true,"layers":[{"text": "Fonds", "expanded": true, "checked": false, "hidden": false, "disabled": false, "leaf": false, "nodeType" : "gx_layercontainer", "nodeGroup": "6" }, {"text": "BD Ortho", "expanded": false, "checked": false, "hidden": false, "disabled": false, "leaf": true, "nodeType" : "gx_layer", "layer": "ORTHOIMAGERY.ORTHOPHOTOS" }, {"text": "Results", "expanded": false, "checked": false, "hidden": false, "disabled": true, "leaf": true, "nodeType" : "gx_layer", "layer": "result_location" }, {"text": "Localisations", "expanded": false, "checked": true, "hidden": false, "disabled": false, "leaf": true, "nodeType" : "gx_layer", "layer": "all_locations" }, {"text": "Détails sur les localisations", "expanded": false, "checked": false, "hidden": true, "disabled": false, "leaf": true, "nodeType" : "gx_layer", "layer": "plot_location_detail" }]}'

$ws.Activate()
$ws.Range("M31").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 6